# "corrected incorrect lists ;)"
# The workbook under test ("Tabelle1") is a fill-in-the-blanks style list of
# two-letter syllables in columns B:G. The author fixed a batch of cells
# that had been shuffled/typo'd, and also moved the active selection.
# Re-apply every corrected cell value (verbatim two-letter strings that are
# stored as shared strings in the original file) and restore the reported
# selection/active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("F2").Value  = "ko"

$ws.Range("B4").Value  = "fi"
$ws.Range("C4").Value  = "fu"
$ws.Range("D4").Value  = "lo"
$ws.Range("E4").Value  = "fi"
$ws.Range("F4").Value  = "lo"

$ws.Range("C5").Value  = "pe"
$ws.Range("D5").Value  = "ba"
$ws.Range("E5").Value  = "lo"
$ws.Range("F5").Value  = "ba"

$ws.Range("D7").Value  = "se"
$ws.Range("E7").Value  = "do"

$ws.Range("F8").Value  = "to"

$ws.Range("B9").Value  = "lo"
$ws.Range("C9").Value  = "se"

$ws.Range("B10").Value = "mi"
$ws.Range("C10").Value = "la"
$ws.Range("E10").Value = "se"
$ws.Range("F10").Value = "fi"

$ws.Range("D12").Value = "fu"
$ws.Range("E12").Value = "ba"

$ws.Range("B16").Value = "se"
$ws.Range("C16").Value = "to"
$ws.Range("D16").Value = "fu"
$ws.Range("E16").Value = "pe"
$ws.Range("F16").Value = "pe"
$ws.Range("G16").Value = "to"

$ws.Range("E19").Value = "se"
$ws.Range("F19").Value = "fi"
$ws.Range("G19").Value = "ba"

$ws.Range("B20").Value = "lo"
$ws.Range("C20").Value = "fi"
$ws.Range("D20").Value = "ba"
$ws.Range("F20").Value = "lo"
$ws.Range("G20").Value = "ba"

$ws.Range("D23").Value = "fu"
$ws.Range("F23").Value = "bo"

$ws.Range("B24").Value = "fu"
$ws.Range("E24").Value = "pe"
$ws.Range("F24").Value = "fu"
$ws.Range("G24").Value = "to"

$ws.Range("B26").Value = "ka"
$ws.Range("C26").Value = "fi"
$ws.Range("D26").Value = "di"
$ws.Range("G26").Value = "ba"

$ws.Range("B27").Value = "pe"
$ws.Range("C27").Value = "li"
$ws.Range("E27").Value = "ba"
$ws.Range("F27").Value = "pe"

$ws.Range("B28").Value = "fu"
$ws.Range("E28").Value = "pe"
$ws.Range("F28").Value = "fu"

$ws.Range("F30").Value = "ni"

$ws.Range("D31").Value = "ki"
$ws.Range("E31").Value = "fi"

$ws.Range("D32").Value = "ba"
$ws.Range("F32").Value = "ti"

$ws.Range("C33").Value = "to"
$ws.Range("F33").Value = "ba"
$ws.Range("G33").Value = "lo"

$ws.Range("C34").Value = "fi"
$ws.Range("F34").Value = "lo"

$ws.Range("D36").Value = "fu"
$ws.Range("F36").Value = "ba"

$ws.Range("C37").Value = "lo"
$ws.Range("E37").Value = "pe"
$ws.Range("F37").Value = "fi"

$ws.Range("C38").Value = "lo"
$ws.Range("D38").Value = "se"

$ws.Range("B40").Value = "to"
$ws.Range("F40").Value = "pe"

# Restore the reported view/selection: scrolled so row 18 is at the top,
# with C28 as the active cell.
$ws.Activate()
[void]$ws.Range("C28").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
